$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("S4:X53").ClearContents()
$ws.Range("S4").Select()
$excel.ActiveCell.End(-4121).Select()
Write-Host "ActiveCell:" $excel.ActiveCell.Address()
